$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 3 DataSource is regressed from the "i-" preprod host to the plain
# preprod host, plus updated Documento/PAS numbers, to do regression
# testing on PreProd.
$ws.Range("A3").Value = "preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("G3").Value = 30633657625
$ws.Range("M3").Value = 304

# Keep the B3 hyperlink target in sync with the new URL (leave the B2
# hyperlink, which still points at the ssurgwsoadev4 host, untouched).
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$3') {
        $hl.Address = "https://preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
    }
}
